$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "קופסה קטנה"
$ws.Range("B3").Value = "קופסה בינונית"
$ws.Range("B4").Value = "שכיבה"
$ws.Range("B5").Value = "קופסת בגדים"
